$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowRange($ws, $rangeAddr, $values) {
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range($rangeAddr).Value = $arr
}

Set-RowRange $ws "B342:AC342" @(6847943, "France Ligue 1", "France Ligue 1", 45263.45833333334, "Monaco", "Montpellier", 2, 0, "H", 1.571, 4, 5.5, 1.571, 4.2, 5.5, -1, 1.99, 1.94, 3, 1.85, 2.05, 0.571, -1, -1, 0.99, -1, -1, 1.05)
Set-RowRange $ws "B343:AC343" @(6847941, "France Ligue 1", "France Ligue 1", 45263.45833333334, "Toulouse", "Lorient", 1, 1, "D", 2.2, 3.3, 3.25, 1.909, 3.3, 4.333, -0.5, 1.975, 1.875, 2.5, 1.975, 1.875, -1, 2.3, -1, -1, 0.875, -1, 0.875)
Set-RowRange $ws "B371:AC371" @(6847975, "France Ligue 1", "France Ligue 1", 45280.70833333334, "Montpellier", "Marseille", 1, 1, "D", 3.4, 3.4, 2.1, 4, 3.6, 1.95, 0.5, 1.9, 1.95, 2.5, 1.925, 1.925, -1, 2.6, -1, 0.8999999999999999, -1, -1, 0.925)
Set-RowRange $ws "B372:AC372" @(6847968, "France Ligue 1", "France Ligue 1", 45280.70833333334, "PSG", "Metz", 3, 1, "H", 1.142, 8, 17, 1.111, 8.5, 23, -2.25, 1.87, 2.06, 3.25, 1.87, 2.03, 0.111, -1, -1, -0.5, 0.53, 0.8700000000000001, -1)
Set-RowRange $ws "B373:AC373" @(6847974, "France Ligue 1", "France Ligue 1", 45280.70833333334, "Reims", "Le Havre", 1, 0, "H", 1.65, 4, 5, 1.727, 3.8, 4.75, -0.75, 1.975, 1.875, 2.5, 1.925, 1.925, 0.7270000000000001, -1, -1, 0.4875, -0.5, -1, 0.925)
Set-RowRange $ws "B374:AC374" @(7579908, "France Ligue 1", "France Ligue 1", 45280.70833333334, "Clermont Foot", "Rennes", 1, 3, "A", 3.25, 3.3, 2.2, 3.1, 3.25, 2.375, 0.25, 1.85, 2.08, 2.25, 2.11, 1.79, -1, -1, 1.375, -1, 1.08, 1.11, -1)

Set-RowRange $ws "B423:G423" @(6848029, "France Ligue 1", "France Ligue 1", 45347.45833333334, "Le Havre", "Reims")
Set-RowRange $ws "K423:AA423" @(3.2, 3.25, 2.25, 3.1, 3.25, 2.375, 0.25, 1.83, 2.07, 2.25, 1.95, 1.95, 0, 0, 0, 0, 0)
Set-RowRange $ws "B425:G425" @(6848027, "France Ligue 1", "France Ligue 1", 45347.45833333334, "Nice", "Clermont Foot")
Set-RowRange $ws "K425:AA425" @(1.45, 4.2, 7.5, 1.4, 4.5, 8, -1.25, 2.07, 1.83, 2.25, 1.83, 2.07, 0, 0, 0, 0, 0)
